$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 2, shifting existing rows down
$ws.Rows.Item(2).Insert()

# Populate the new row with the new survey entry
# (shared-string insertion order matters: archive_id, wave, description, timeframe)
$ws.Cells.Item(2, 1).Value = "ZA7953"
$ws.Cells.Item(2, 2).Value = "'98.2"
$ws.Cells.Item(2, 4).Value = "Standard Eurobarometer 98 (COVID-19 Pandemic)"
$ws.Cells.Item(2, 3).Value = "January-February 2023"

# Update the active selection to match the post-edit state
$ws.Range("C3").Select()
